$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 29200
$ws.Range("B6").Value = "OWN"
$ws.Range("C6").Value = "High"
$ws.Range("D6").Value = 800000
$ws.Range("C8").Value = "High"
$ws.Range("D8").Value = 470004
$ws.Range("D9").Value = 18000
$ws.Range("D10").Value = 30000
$ws.Range("C12").Value = "Medium"
$ws.Range("C13").Value = "Medium"
$ws.Range("B16").Value = "OWN"
$ws.Range("C16").Value = "Medium"
$ws.Range("D16").Value = 50000
$ws.Range("D18").Value = 38000
$ws.Range("D19").Value = 61000
$ws.Range("D20").Value = 34000
$ws.Range("D21").Value = 17385
$ws.Range("D22").Value = 23370
$ws.Range("C23").Value = "Low"
$ws.Range("D23").Value = 10500
$ws.Range("C24").Value = "Medium"
$ws.Range("C25").Value = "Medium"
$ws.Range("C26").Value = "Medium"
$ws.Range("C27").Value = "Medium"
$ws.Range("D28").Value = 30000
$ws.Range("D29").Value = 30004
$ws.Range("C30").Value = "High"
$ws.Range("C33").Value = "Medium"
$ws.Range("D37").Value = 31000
$ws.Range("D38").Value = 35596
$ws.Range("D39").Value = 35000
$ws.Range("D42").Value = 38004
$ws.Range("D43").Value = 32300
$ws.Range("D44").Value = 35000
$ws.Range("D45").Value = 35000
$ws.Range("C46").Value = "Medium"
$ws.Range("C47").Value = "Medium"
$ws.Range("C48").Value = "High"
$ws.Range("D50").Value = 30000
$ws.Range("D51").Value = 30000
$ws.Range("D52").Value = 20000
$ws.Range("D53").Value = 70000
$ws.Range("D54").Value = 20000
$ws.Range("C55").Value = "Medium"
$ws.Range("D56").Value = 26000
$ws.Range("C57").Value = "Medium"
$ws.Range("C58").Value = "Medium"
$ws.Range("C60").Value = "Medium"
$ws.Range("D63").Value = 39000
$ws.Range("D64").Value = 39000
$ws.Range("C65").Value = "Medium"
$ws.Range("C66").Value = "Medium"
$ws.Range("C67").Value = "Medium"
$ws.Range("D69").Value = 34544
$ws.Range("D70").Value = 25000
$ws.Range("D71").Value = 21000
$ws.Range("D72").Value = 24000
$ws.Range("D73").Value = 35000
$ws.Range("D75").Value = 30000
$ws.Range("C76").Value = "High"
$ws.Range("C77").Value = "Medium"
$ws.Range("C78").Value = "Medium"
$ws.Range("C79").Value = "Medium"
$ws.Range("C80").Value = "Medium"
$ws.Range("C81").Value = "Medium"
$ws.Range("C82").Value = "High"
$ws.Range("D83").Value = 38000
$ws.Range("C84").Value = "High"
$ws.Range("D84").Value = 332000
$ws.Range("C85").Value = "High"
$ws.Range("D85").Value = 153000
$ws.Range("B86").Value = "OWN"
$ws.Range("C87").Value = "Medium"
$ws.Range("C89").Value = "High"
$ws.Range("D89").Value = 242000
$ws.Range("D90").Value = 28365
$ws.Range("D91").Value = 22000
$ws.Range("D92").Value = 22000
$ws.Range("D93").Value = 33680
$ws.Range("C94").Value = "High"
$ws.Range("C95").Value = "Medium"
$ws.Range("C97").Value = "Medium"
$ws.Range("C98").Value = "Medium"
$ws.Range("D99").Value = 35000
$ws.Range("B100").Value = "OWN"
$ws.Range("C100").Value = "High"

# Update the saved selection/view state to match the author's last position
$ws.Range("D99").Select()
